$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1397.5883
$ws.Range("J17").Value = 1413.0625
$ws.Range("L17").Value = 4239.1875
$ws.Range("N17").Value = -4575.1875
$ws.Range("H80").Value = 1588
$ws.Range("I80").Value = 596.8333
$ws.Range("J80").Value = 2331.375
$ws.Range("K80").Value = 1790.4999
$ws.Range("L80").Value = 6994.125
$ws.Range("M80").Value = -792.4999
$ws.Range("N80").Value = -8990.125
$ws.Range("H83").Value = 1588
$ws.Range("I83").Value = 596.8333
$ws.Range("J83").Value = 2331.375
$ws.Range("K83").Value = 5371.4997
$ws.Range("L83").Value = 20982.375
$ws.Range("M83").Value = -379.4997000000003
$ws.Range("N83").Value = -30966.375
$ws.Range("H86").Value = 2550.5625
$ws.Range("I86").Value = 3067.8
$ws.Range("K86").Value = 3067.8
$ws.Range("M86").Value = -1944.8
$ws.Range("H89").Value = 2550.5625
$ws.Range("I89").Value = 3067.8
$ws.Range("K89").Value = 15339
$ws.Range("M89").Value = -9723
$ws.Range("H125").Value = 7248967.5
$ws.Range("I125").Value = 1616.3334
$ws.Range("J125").Value = 8336070.5
$ws.Range("K125").Value = 14547.0006
$ws.Range("L125").Value = 75024634.5
$ws.Range("M125").Value = -12087.0006
$ws.Range("N125").Value = -75029554.5
$ws.Range("H132").Value = 20411054
$ws.Range("I132").Value = 21742124
$ws.Range("K132").Value = 65226372
$ws.Range("M132").Value = -65223842
$ws.Range("H137").Value = 53072.06
$ws.Range("I137").Value = 138988.53
$ws.Range("K137").Value = 416965.59
$ws.Range("M137").Value = -414415.59
$ws.Range("H138").Value = 2641.414
$ws.Range("I138").Value = 1532.3636
$ws.Range("K138").Value = 4597.0908
$ws.Range("M138").Value = 542.9092000000001
$ws.Range("H141").Value = 1632.1305
$ws.Range("I141").Value = 1456.8125
$ws.Range("J141").Value = 2032.8572
$ws.Range("K141").Value = 4370.4375
$ws.Range("L141").Value = 6098.571599999999
$ws.Range("M141").Value = 809.5625
$ws.Range("N141").Value = -16458.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2689.8572
$ws.Range("I63").Value = 2596.7693
$ws.Range("K63").Value = 2596.7693
$ws.Range("M63").Value = -1910.7693
$ws.Range("H66").Value = 2689.8572
$ws.Range("I66").Value = 2596.7693
$ws.Range("K66").Value = 12983.8465
$ws.Range("M66").Value = -9551.8465
$ws.Range("H74").Value = 50134.41
$ws.Range("I74").Value = 3624.5227
$ws.Range("K74").Value = 3624.5227
$ws.Range("M74").Value = -2750.5227
$ws.Range("H77").Value = 50134.41
$ws.Range("I77").Value = 3624.5227
$ws.Range("K77").Value = 18122.6135
$ws.Range("M77").Value = -13754.6135
$ws.Range("H110").Value = 1992278.9
$ws.Range("I110").Value = 2533264.2
$ws.Range("K110").Value = 2533264.2
$ws.Range("M110").Value = -2531219.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3843.2058
$ws.Range("I134").Value = 1562
$ws.Range("J134").Value = 6732.7334
$ws.Range("K134").Value = 4686
$ws.Range("L134").Value = 20198.2002
$ws.Range("M134").Value = -2151
$ws.Range("N134").Value = -25268.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H31").Value = 21801.8
$ws.Range("I31").Value = 1044.4062
$ws.Range("J31").Value = 72896.92
$ws.Range("K31").Value = 1044.4062
$ws.Range("L31").Value = 72896.92
$ws.Range("M31").Value = -749.4061999999999
$ws.Range("N31").Value = -73486.92
$ws.Range("H34").Value = 21801.8
$ws.Range("I34").Value = 1044.4062
$ws.Range("J34").Value = 72896.92
$ws.Range("K34").Value = 1044.4062
$ws.Range("L34").Value = 72896.92
$ws.Range("M34").Value = -842.4061999999999
$ws.Range("N34").Value = -73300.92
$ws.Range("H132").Value = 71888.25999999999
$ws.Range("I132").Value = 68363.92999999999
$ws.Range("J132").Value = 76293.664
$ws.Range("K132").Value = 205091.79
$ws.Range("L132").Value = 228880.992
$ws.Range("M132").Value = -202561.79
$ws.Range("N132").Value = -233940.992
$ws.Range("H134").Value = 4344
$ws.Range("I134").Value = 3897.8667
$ws.Range("J134").Value = 5087.5557
$ws.Range("K134").Value = 11693.6001
$ws.Range("L134").Value = 15262.6671
$ws.Range("M134").Value = -9158.6001
$ws.Range("N134").Value = -20332.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 192.36842
$ws.Range("J107").Value = 283
$ws.Range("L107").Value = 849
$ws.Range("N107").Value = -4689
$ws.Range("H113").Value = 5029.7856
$ws.Range("I113").Value = 9980
$ws.Range("K113").Value = 29940
$ws.Range("M113").Value = -27770
$ws.Range("H138").Value = 4699.6665
$ws.Range("I138").Value = 3785
$ws.Range("J138").Value = 5500
$ws.Range("K138").Value = 11355
$ws.Range("L138").Value = 16500
$ws.Range("M138").Value = -6215
$ws.Range("N138").Value = -26780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1527763
$ws.Range("I80").Value = 3050887.2
$ws.Range("J80").Value = 4638.875
$ws.Range("K80").Value = 3050887.2
$ws.Range("L80").Value = 4638.875
$ws.Range("M80").Value = -3049889.2
$ws.Range("N80").Value = -6634.875
$ws.Range("H83").Value = 1527763
$ws.Range("I83").Value = 3050887.2
$ws.Range("J83").Value = 4638.875
$ws.Range("K83").Value = 15254436
$ws.Range("L83").Value = 23194.375
$ws.Range("M83").Value = -15249444
$ws.Range("N83").Value = -33178.375
$ws.Range("H113").Value = 33335172
$ws.Range("I113").Value = 33335172
$ws.Range("K113").Value = 33335172
$ws.Range("M113").Value = -33333002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6333
$ws.Range("I7").Value = 3666.3333
$ws.Range("J7").Value = 7666.3335
$ws.Range("K7").Value = 3666.3333
$ws.Range("L7").Value = 7666.3335
$ws.Range("N7").Value = -7890.3335
$ws.Range("H40").Value = 8167.7
$ws.Range("I40").Value = 6739.4287
$ws.Range("K40").Value = 6739.4287
$ws.Range("M40").Value = -6603.4287
$ws.Range("H46").Value = 5879.7
$ws.Range("J46").Value = 6074.625
$ws.Range("L46").Value = 6074.625
$ws.Range("N46").Value = -6450.625
$ws.Range("H61").Value = 13894322
$ws.Range("J61").Value = 3937.5
$ws.Range("L61").Value = 3937.5
$ws.Range("N61").Value = -4341.5
$ws.Range("H113").Value = 13894322
$ws.Range("J113").Value = 3937.5
$ws.Range("L113").Value = 3937.5
$ws.Range("N113").Value = -8277.5
$ws.Range("H126").Value = 6333
$ws.Range("I126").Value = 3666.3333
$ws.Range("J126").Value = 7666.3335
$ws.Range("K126").Value = 10998.9999
$ws.Range("L126").Value = 22999.0005
$ws.Range("N126").Value = -27939.0005
$ws.Range("H136").Value = 66920.34
$ws.Range("I136").Value = 83542.08
$ws.Range("K136").Value = 250626.24
$ws.Range("M136").Value = -248076.24

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8303.071
$ws.Range("J62").Value = 9647.695
$ws.Range("L62").Value = 9647.695
$ws.Range("N62").Value = -10895.695
$ws.Range("H65").Value = 8303.071
$ws.Range("J65").Value = 9647.695
$ws.Range("L65").Value = 48238.475
$ws.Range("N65").Value = -54478.475
$ws.Range("H122").Value = 2231.9
$ws.Range("I122").Value = 1532.909
$ws.Range("K122").Value = 4598.727000000001
$ws.Range("M122").Value = -2148.727000000001
$ws.Range("H126").Value = 2595.7144
$ws.Range("I126").Value = 3577.8572
$ws.Range("J126").Value = 1613.5714
$ws.Range("K126").Value = 10733.5716
$ws.Range("L126").Value = 4840.7142
$ws.Range("M126").Value = -8263.571599999999
$ws.Range("N126").Value = -9780.7142
$ws.Range("H132").Value = 22976106
$ws.Range("I132").Value = 28572482
$ws.Range("K132").Value = 85717446
$ws.Range("M132").Value = -85714916
